$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.770.26'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').Value = '2.651.05'
$ws.Range('E3').Value = '  +2.06%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '536.85'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = '145.70'
$ws.Range('E6').Value = '  +3.38%  '
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '0.574'
$ws.Range('E8').Value = '  +1.10%  '
$ws.Range('D9').Value = '2.666.44'
$ws.Range('E9').Value = '  +2.09%  '
$ws.Range('D10').Value = '6.69'
$ws.Range('E10').Value = '  +3.37%  '
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').Value = '0.338'
$ws.Range('E12').Value = '  +0.94%  '
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').Value = '3.117.17'
$ws.Range('E14').Value = '  +1.88%  '
$ws.Range('D15').Value = '59.667.25'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('D16').Value = '21.19'
$ws.Range('E16').Value = '  +3.20%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = '0.0000135'
$ws.Range('E17').Value = '  +1.17%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.631.08'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('D19').Value = '344.94'
$ws.Range('E19').Value = '  -0.40%  '
$ws.Range('D20').Value = '4.42'
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('D21').Value = '10.25'
$ws.Range('E21').Value = '  +1.00%  '
$ws.Range('D22').Value = '6.36'
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').Value = '66.64'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').Value = '0.416'
$ws.Range('E25').Value = '  +2.21%  '
$ws.Range('E26').Value = '  -1.53%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  +0.13%  '
$ws.Range('D28').Value = '7.30'
$ws.Range('E28').Value = '  +1.12%  '
$ws.Range('D29').Value = '0.0₃0754'
$ws.Range('E29').Value = '  +2.08%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = '1.66'
$ws.Range('E31').Value = '  +1.41%  '
$ws.Range('D32').Value = '5.85'
$ws.Range('E32').Value = '  +0.39%  '
$ws.Range('D33').Value = '19.02'
$ws.Range('E33').Value = '  +0.94%  '
$ws.Range('D34').Value = '149.89'
$ws.Range('E34').Value = '  +0.49%  '
$ws.Range('D35').Value = '4.03'
$ws.Range('E35').Value = '  +0.86%  '
$ws.Range('E36').Value = '  +2.26%  '
$ws.Range('D37').Value = '0.844'
$ws.Range('E37').Value = '  -0.48%  '
$ws.Range('E38').Value = '  -1.31%  '
$ws.Range('D39').Value = '0.827'
$ws.Range('E39').Value = '  -0.95%  '
$ws.Range('D40').Value = '291.90'
$ws.Range('E40').Value = '  +5.10%  '
$ws.Range('D41').Value = '3.61'
$ws.Range('E41').Value = '  +1.97%  '
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.08%  '
$ws.Range('D43').Value = '0.606'
$ws.Range('E43').Value = '  +1.16%  '
$ws.Range('D44').Value = '0.0544'
$ws.Range('E44').Value = '  +4.27%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '19.46'
$ws.Range('E45').Value = '  +4.94%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value = '10.74'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('E47').Value = '  -1.19%  '
$ws.Range('D48').Value = '1.976.15'
$ws.Range('E48').Value = '  +1.57%  '
$ws.Range('E49').Value = '  +1.64%  '
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('D51').Value = '18.40'
$ws.Range('E51').Value = '  +0.35%  '
